$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H48").Value = 57765.5
$ws.Range("I48").Value = 1017
$ws.Range("K48").Value = 3051
$ws.Range("M48").Value = -2759
$ws.Range("H56").Value = 57765.5
$ws.Range("I56").Value = 1017
$ws.Range("K56").Value = 3051
$ws.Range("M56").Value = -2517
$ws.Range("H112").Value = 5064.7617
$ws.Range("I112").Value = 420
$ws.Range("J112").Value = 5553.684
$ws.Range("K112").Value = 1260
$ws.Range("L112").Value = 16661.052
$ws.Range("M112").Value = -152
$ws.Range("N112").Value = -18877.052
$ws.Range("H113").Value = 2375.9412
$ws.Range("I113").Value = 2263.125
$ws.Range("J113").Value = 2476.2222
$ws.Range("K113").Value = 2263.125
$ws.Range("L113").Value = 2476.2222
$ws.Range("M113").Value = 990.875
$ws.Range("N113").Value = -8984.2222
$ws.Range("H132").Value = 1977.5902
$ws.Range("I132").Value = 1974.2354
$ws.Range("J132").Value = 1994.7
$ws.Range("K132").Value = 5922.706200000001
$ws.Range("L132").Value = 5984.1
$ws.Range("M132").Value = -3392.706200000001
$ws.Range("N132").Value = -11044.1
$ws.Range("H137").Value = 2977840
$ws.Range("I137").Value = 9260471
$ws.Range("K137").Value = 27781413
$ws.Range("M137").Value = -27778863
$ws.Range("H138").Value = 3331.394
$ws.Range("I138").Value = 2234.0386
$ws.Range("J138").Value = 3722.233
$ws.Range("K138").Value = 6702.1158
$ws.Range("L138").Value = 11166.699
$ws.Range("M138").Value = -1562.1158
$ws.Range("N138").Value = -21446.699
$ws.Range("H140").Value = 74048.5
$ws.Range("J140").Value = 74048.5
$ws.Range("L140").Value = 74048.5
$ws.Range("N140").Value = -84408.5
$ws.Range("H141").Value = 3417.96
$ws.Range("I141").Value = 1586.8889
$ws.Range("J141").Value = 8126.4287
$ws.Range("K141").Value = 4760.6667
$ws.Range("L141").Value = 24379.2861
$ws.Range("M141").Value = 419.3333000000002
$ws.Range("N141").Value = -34739.2861

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 14495818
$ws.Range("I61").Value = 22224320
$ws.Range("K61").Value = 22224320
$ws.Range("M61").Value = -22224108
$ws.Range("H63").Value = 3942.348
$ws.Range("I63").Value = 2838.9
$ws.Range("K63").Value = 2838.9
$ws.Range("M63").Value = -2152.9
$ws.Range("H66").Value = 3942.348
$ws.Range("I66").Value = 2838.9
$ws.Range("K66").Value = 14194.5
$ws.Range("M66").Value = -10762.5
$ws.Range("H74").Value = 13891615
$ws.Range("I74").Value = 2137.6316
$ws.Range("J74").Value = 29415148
$ws.Range("K74").Value = 2137.6316
$ws.Range("L74").Value = 29415148
$ws.Range("M74").Value = -1263.6316
$ws.Range("N74").Value = -29416896
$ws.Range("H76").Value = 72729.336
$ws.Range("J76").Value = 72729.336
$ws.Range("L76").Value = 72729.336
$ws.Range("N76").Value = -73405.336
$ws.Range("H77").Value = 13891615
$ws.Range("I77").Value = 2137.6316
$ws.Range("J77").Value = 29415148
$ws.Range("K77").Value = 10688.158
$ws.Range("L77").Value = 147075740
$ws.Range("M77").Value = -6320.158000000001
$ws.Range("N77").Value = -147084476
$ws.Range("H79").Value = 72729.336
$ws.Range("J79").Value = 72729.336
$ws.Range("L79").Value = 72729.336
$ws.Range("N79").Value = -75069.336
$ws.Range("H132").Value = 1835224.4
$ws.Range("I132").Value = 3104.08
$ws.Range("J132").Value = 4529519
$ws.Range("K132").Value = 9312.24
$ws.Range("L132").Value = 13588557
$ws.Range("M132").Value = -6782.24
$ws.Range("N132").Value = -13593617
$ws.Range("H136").Value = 14495818
$ws.Range("I136").Value = 22224320
$ws.Range("K136").Value = 66672960
$ws.Range("M136").Value = -66670410

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 470.9375
$ws.Range("I22").Value = 336.625
$ws.Range("J22").Value = 605.25
$ws.Range("K22").Value = 336.625
$ws.Range("L22").Value = 605.25
$ws.Range("M22").Value = 13.375
$ws.Range("N22").Value = -1305.25
$ws.Range("H31").Value = 5251.0366
$ws.Range("I31").Value = 1781.7142
$ws.Range("J31").Value = 6445.3936
$ws.Range("K31").Value = 1781.7142
$ws.Range("L31").Value = 6445.3936
$ws.Range("M31").Value = -1486.7142
$ws.Range("N31").Value = -7035.3936
$ws.Range("H34").Value = 5251.0366
$ws.Range("I34").Value = 1781.7142
$ws.Range("J34").Value = 6445.3936
$ws.Range("K34").Value = 1781.7142
$ws.Range("L34").Value = 6445.3936
$ws.Range("M34").Value = -1579.7142
$ws.Range("N34").Value = -6849.3936
$ws.Range("H58").Value = 1241.5
$ws.Range("I58").Value = 937.5
$ws.Range("J58").Value = 2153.5
$ws.Range("K58").Value = 937.5
$ws.Range("L58").Value = 2153.5
$ws.Range("M58").Value = -734.5
$ws.Range("N58").Value = -2559.5
$ws.Range("H86").Value = 2361.9524
$ws.Range("I86").Value = 2322.389
$ws.Range("J86").Value = 2599.3333
$ws.Range("K86").Value = 2322.389
$ws.Range("L86").Value = 2599.3333
$ws.Range("M86").Value = -1199.389
$ws.Range("N86").Value = -4845.3333
$ws.Range("H89").Value = 2361.9524
$ws.Range("I89").Value = 2322.389
$ws.Range("J89").Value = 2599.3333
$ws.Range("K89").Value = 11611.945
$ws.Range("L89").Value = 12996.6665
$ws.Range("M89").Value = -5995.945
$ws.Range("N89").Value = -24228.6665
$ws.Range("H99").Value = 2461.0789
$ws.Range("I99").Value = 2312.2
$ws.Range("J99").Value = 2514.25
$ws.Range("K99").Value = 2312.2
$ws.Range("L99").Value = 2514.25
$ws.Range("M99").Value = -814.1999999999998
$ws.Range("N99").Value = -5510.25
$ws.Range("H126").Value = 2461.0789
$ws.Range("I126").Value = 2312.2
$ws.Range("J126").Value = 2514.25
$ws.Range("K126").Value = 6936.599999999999
$ws.Range("L126").Value = 7542.75
$ws.Range("M126").Value = -4466.599999999999
$ws.Range("N126").Value = -12482.75
$ws.Range("H132").Value = 16261999
$ws.Range("I132").Value = 19232338
$ws.Range("K132").Value = 57697014
$ws.Range("M132").Value = -57694484
$ws.Range("H136").Value = 1241.5
$ws.Range("I136").Value = 937.5
$ws.Range("J136").Value = 2153.5
$ws.Range("K136").Value = 2812.5
$ws.Range("L136").Value = 6460.5
$ws.Range("M136").Value = -262.5
$ws.Range("N136").Value = -11560.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 7470.5835
$ws.Range("I3").Value = 4847.143
$ws.Range("K3").Value = 14541.429
$ws.Range("M3").Value = -14429.429
$ws.Range("H113").Value = 592.5238000000001
$ws.Range("I113").Value = 662.8570999999999
$ws.Range("J113").Value = 557.3570999999999
$ws.Range("K113").Value = 1988.5713
$ws.Range("L113").Value = 1672.0713
$ws.Range("M113").Value = 181.4287000000002
$ws.Range("N113").Value = -6012.0713
$ws.Range("H123").Value = 1583.8334
$ws.Range("I123").Value = 1015
$ws.Range("J123").Value = 2721.5
$ws.Range("K123").Value = 3045
$ws.Range("L123").Value = 8164.5
$ws.Range("M123").Value = -595
$ws.Range("N123").Value = -13064.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5364.1304
$ws.Range("I70").Value = 5226.3887
$ws.Range("J70").Value = 5860
$ws.Range("K70").Value = 5226.3887
$ws.Range("L70").Value = 5860
$ws.Range("M70").Value = -4956.3887
$ws.Range("N70").Value = -6400
$ws.Range("H73").Value = 5364.1304
$ws.Range("I73").Value = 5226.3887
$ws.Range("J73").Value = 5860
$ws.Range("K73").Value = 5226.3887
$ws.Range("L73").Value = 5860
$ws.Range("M73").Value = -4290.3887
$ws.Range("N73").Value = -7732
$ws.Range("H113").Value = 84676.664
$ws.Range("I113").Value = 112146.664
$ws.Range("J113").Value = 2266.6667
$ws.Range("K113").Value = 112146.664
$ws.Range("L113").Value = 2266.6667
$ws.Range("M113").Value = -109976.664
$ws.Range("N113").Value = -6606.6667
$ws.Range("H132").Value = 27783242
$ws.Range("I132").Value = 55563356
$ws.Range("J132").Value = 3125.3333
$ws.Range("K132").Value = 166690068
$ws.Range("L132").Value = 9375.999899999999
$ws.Range("M132").Value = -166687538
$ws.Range("N132").Value = -14435.9999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 4268.3335
$ws.Range("I40").Value = 3782
$ws.Range("K40").Value = 3782
$ws.Range("M40").Value = -3646
$ws.Range("H88").Value = 40000
$ws.Range("I88").Value = 40000
$ws.Range("K88").Value = 40000
$ws.Range("M88").Value = -39572
$ws.Range("H91").Value = 40000
$ws.Range("I91").Value = 40000
$ws.Range("K91").Value = 40000
$ws.Range("M91").Value = -38518
$ws.Range("H133").Value = 50560.855
$ws.Range("J133").Value = 50560.855
$ws.Range("L133").Value = 50560.855
$ws.Range("N133").Value = -55620.855

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 3321
$ws.Range("I81").Value = 3466
$ws.Range("J81").Value = 2849.75
$ws.Range("K81").Value = 6932
$ws.Range("L81").Value = 5699.5
$ws.Range("M81").Value = -5871
$ws.Range("N81").Value = -7821.5
$ws.Range("H84").Value = 3321
$ws.Range("I84").Value = 3466
$ws.Range("J84").Value = 2849.75
$ws.Range("K84").Value = 34660
$ws.Range("L84").Value = 28497.5
$ws.Range("M84").Value = -29356
$ws.Range("N84").Value = -39105.5
$ws.Range("H132").Value = 5305209.5
$ws.Range("I132").Value = 2135.0278
$ws.Range("J132").Value = 15353141
$ws.Range("K132").Value = 6405.0834
$ws.Range("L132").Value = 46059423
$ws.Range("M132").Value = -3875.0834
$ws.Range("N132").Value = -46064483
